# Regenerate save_data to use K instead of Strike#, recalculated std/mean,
# and write the newly calculated s_vals into column G (K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K (column G) value, per recalculated s_vals
$newValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    12 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
